$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.497.56"
$ws.Range("E2").Value = "  -3.53%  "
$ws.Range("D3").Value = "1.959.49"
$ws.Range("E3").Value = "  -1.75%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.012"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.31%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "321.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.010"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.38%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4769"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4068"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.25%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.43"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08465"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.97%  "
$ws.Range("E11").Value = "  -4.33%  "
$ws.Range("E12").Value = "  -4.15%  "
$ws.Range("D13").Value = "1.969.04"
$ws.Range("E13").Value = "  -1.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.619"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.85%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.169"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.74%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.012"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "89.40"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.38%  "
$ws.Range("E18").Value = "  -2.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06614"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.74%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.62"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.97%  "
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.830"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.88%  "
$ws.Range("D23").Value = "28.518.45"
$ws.Range("E23").Value = "  -3.51%  "
$ws.Range("E24").Value = "  -2.63%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.290"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.47%  "
$ws.Range("D26").Value = "2.206.62"
$ws.Range("E26").Value = "  -0.66%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "155.21"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.23"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.88%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.937"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.93%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.162"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "123.75"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.60%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9816"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.87%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09616"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.57%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.445"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.42%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.597"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.10%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.664"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.81%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02336"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.29%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.829"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06231"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.81%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.260"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.79%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6230"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.16"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.17%  "
$ws.Range("E43").Value = "  -0.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1922"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.22%  "
$ws.Range("E45").Value = "  +2.58%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5969"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "13.01"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.060"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.40%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.403"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.81%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00000000328"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06836"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.59%  "

Write-Output "Applied 97 cell updates"